$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis) " heading paragraph.
$findRange = $d.Content
[void]$findRange.Find.Execute("Docente(s) Responsável(eis) ")
$docentePara = $findRange.Paragraphs(1)

# Insert a brand-new (empty) paragraph right after it, then grab a
# handle to that freshly created paragraph.
[void]$docentePara.Range.InsertParagraphAfter()
$newPara = $docentePara.Next()

# Populate it with the ListBullet-styled content: two names separated by
# a manual line break, each its own run (matching the target markup).
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>5817692 - Katia Cristiane Gandolpho Candioto</w:t><w:br/></w:r><w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t></w:r></w:p>'
[void]$newPara.Range.InsertXML($xml)
